# New crime data collected — weekly CompStat refresh (109th Precinct)
# Updates: report header (volume/date range) + weekly crime-stat grid (rows 15-21,
# 22, 24-28, 33), including a few cells that flip between numeric and the
# sheet's "N/A" text markers ("0" / "***.*").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings -------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/24/2024  Through  6/30/2024"

# --- Row 15 (Rape) ---------------------------------------------------------
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 100
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -15.789473684210
$ws.Range("M15").Value = 77.777777777777
$ws.Range("N15").Value = 14.285714285714

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -4.545454545454
$ws.Range("I16").Value = 175
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = 2.941176470588
$ws.Range("L16").Value = 35.658914728682
$ws.Range("M16").Value = 69.902912621359
$ws.Range("N16").Value = -69.405594405594

# --- Row 17 (Fel. Assault) --------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = 37.837837837837
$ws.Range("I17").Value = 237
$ws.Range("J17").Value = 237
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 68.085106382978
$ws.Range("M17").Value = 137
$ws.Range("N17").Value = 33.146067415730

# --- Row 18 (Burglary) ------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -72.727272727272
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -31.25
$ws.Range("I18").Value = 189
$ws.Range("J18").Value = 288
$ws.Range("K18").Value = -34.375
$ws.Range("L18").Value = -16
$ws.Range("M18").Value = -15.246636771300
$ws.Range("N18").Value = -83.996613039796

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 28
$ws.Range("E19").Value = -3.448275862068
$ws.Range("F19").Value = 96
$ws.Range("G19").Value = 98
$ws.Range("H19").Value = -2.040816326530
$ws.Range("I19").Value = 657
$ws.Range("J19").Value = 709
$ws.Range("K19").Value = -7.334273624823
$ws.Range("L19").Value = -10.854816824966
$ws.Range("M19").Value = 121.212121212121
$ws.Range("N19").Value = 0.152439024390

# --- Row 20 (G.L.A.) --------------------------------------------------------
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = -57.894736842105
$ws.Range("F20").Value = 41
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = -4.651162790697
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = -0.4
$ws.Range("L20").Value = 139.423076923077
$ws.Range("M20").Value = 79.136690647482
$ws.Range("N20").Value = -87.308868501529

# --- Row 21 (TOTAL) ---------------------------------------------------------
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 74
$ws.Range("E21").Value = -24.324324324324
$ws.Range("F21").Value = 237
$ws.Range("G21").Value = 236
$ws.Range("H21").Value = 0.423728813559
$ws.Range("I21").Value = 1525
$ws.Range("J21").Value = 1674
$ws.Range("K21").Value = -8.900836320191
$ws.Range("L21").Value = 12.546125461254
$ws.Range("M21").Value = 74.485125858123
$ws.Range("N21").Value = -66.615586690017

# --- Row 22 (Transit) -------------------------------------------------------
# C22 flips from the text marker "0" to a real numeric value this week.
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -40
$ws.Range("L22").Value = -25
$ws.Range("M22").Value = 300

# --- Row 24 (Petit Larceny) -------------------------------------------------
$ws.Range("C24").Value = 57
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = 26.666666666666
$ws.Range("F24").Value = 195
$ws.Range("G24").Value = 192
$ws.Range("H24").Value = 1.5625
$ws.Range("I24").Value = 1363
$ws.Range("J24").Value = 1372
$ws.Range("K24").Value = -0.655976676384
$ws.Range("L24").Value = -3.264726756564
$ws.Range("M24").Value = 83.445491251682

# --- Row 25 (Retail Theft) --------------------------------------------------
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 33
$ws.Range("E25").Value = 3.030303030303
$ws.Range("F25").Value = 123
$ws.Range("G25").Value = 120
$ws.Range("H25").Value = 2.5
$ws.Range("I25").Value = 882
$ws.Range("J25").Value = 733
$ws.Range("K25").Value = 20.327421555252
$ws.Range("L25").Value = 22.330097087378

# --- Row 26 (Misd. Assault) -------------------------------------------------
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 84
$ws.Range("H26").Value = 35.483870967741
$ws.Range("I26").Value = 446
$ws.Range("J26").Value = 405
$ws.Range("K26").Value = 10.123456790123
$ws.Range("L26").Value = 31.952662721893
$ws.Range("M26").Value = 25.633802816901

# --- Row 27 (UCR Rape*) ------------------------------------------------------
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = -26.666666666666
$ws.Range("L27").Value = -12

# --- Row 28 (Other Sex Crimes) ----------------------------------------------
# D28/E28 flip from numeric values to the sheet's text markers "0" / "***.*".
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "'***.*"
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 54
$ws.Range("K28").Value = -5.263157894736
$ws.Range("L28").Value = 20

# --- Row 29 (Shooting Vic.) --------------------------------------------------
# D29/E29 flip from numeric values to the text markers "0" / "***.*".
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "'***.*"

# --- Row 30 (Shooting Inc.) --------------------------------------------------
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"

# The assignments above mark D28/E28/D29/E29/D30/E30 as text (quote-prefixed),
# which picks up a throwaway "number stored as text" style. Re-apply the
# workbook's normal "N/A" text style (same one already used by neighboring
# cells such as C29/F29) via a formats-only paste so the cells match the
# look of every other "N/A" cell on the sheet.
$ws.Range("C29").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("F29").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)

# --- Row 33 (Traffic Fatalities) ---------------------------------------------
$ws.Range("J33").Value = 6
$ws.Range("K33").Value = 50
